$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A new weekly record needs to be inserted at the top of the data block
# (row 34), pushing every existing data row down by one. This mirrors the
# "semanal" (weekly) update pattern used throughout this workbook.
$ws.Rows.Item(34).Insert()

# Populate the newly inserted row 34 with the new record's values while
# keeping the same shape as its neighbours (A/B/C/E/F/G/I/N/Q/R repeat the
# same constants used across the whole data block).
$ws.Cells.Item(34, 1).Value = 5
$ws.Cells.Item(34, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(34, 3).Value = "Maule"
$ws.Cells.Item(34, 4).Value = 44914
$ws.Cells.Item(34, 5).Value = 7
$ws.Cells.Item(34, 6).Value = 100112022
$ws.Cells.Item(34, 7).Value = "Arveja Verde"
$ws.Cells.Item(34, 8).Value = "Sin especificar"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 500
$ws.Cells.Item(34, 11).Value = 20000
$ws.Cells.Item(34, 12).Value = 20000
$ws.Cells.Item(34, 13).Value = 20000
$ws.Cells.Item(34, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(34, 15).Value = "Carahue"
$ws.Cells.Item(34, 16).Value = 800
$ws.Cells.Item(34, 17).Value = 25
$ws.Cells.Item(34, 18).Value = "Hortaliza"

# Match the date-number format used by the rest of column D.
$ws.Cells.Item(34, 4).NumberFormat = $ws.Cells.Item(35, 4).NumberFormat
